$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 10 (Objetivos:) content changed in B10/C10
# ---------------------------------------------------------------------
$ws.Range("B10").Value = 'Conduzir os alunos no desenvolvimento de um projeto de conclusão de curso sobre tema específico relacionado à engenharia química.'
$ws.Range("C10").Value = 'Conduzir os alunos no desenvolvimento de um projeto de conclusão de curso sobre tema específico relacionado à engenharia química.'

# ---------------------------------------------------------------------
# 2) Insert a new row at 13 -- shifts old rows 13-21 down to 14-22
#    (old row 12 "Docentes responsaveis:" label now gets its own
#     content row right below it)
# ---------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# The inserted row comes back with a stray styled-but-empty A13 cell; drop it
$ws.Range("A13").Clear()

# Give B13:C13 the normal content-cell formatting (copy from B14:C14) then fill them in
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("B13").Value = '198273 - Domingos Savio Giordani'
$ws.Range("C13").Value = '198273 - Domingos Savio Giordani'

# ---------------------------------------------------------------------
# 3) Refresh the text for rows 14-22 (content shifted down by the insert;
#    several rows also carry new/updated copy per the edit)
# ---------------------------------------------------------------------
# Row 14
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = '1) Metodologia Cientifica. 2) Projeto de Monografia. 3) Métodos de Pesquisa. 4) Normas de um Projeto de Pesquisa. 5) Pesquisa em Bases de Dados Bibliográficos. 6) Organização de Referências Bibliográficas.'
$ws.Range("C14").Value = '1) Metodologia Cientifica. 2) Projeto de Monografia. 3) Métodos de Pesquisa. 4) Normas de um Projeto de Pesquisa. 5) Pesquisa em Bases de Dados Bibliográficos. 6) Organização de Referências Bibliográficas.'

# Row 15
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = '1) Scientific Methodology. 2) Monography Project. 3) Research Methods. 4) Guidelines of a research project. 5) Research in Bibliographic Databases. 6) Organization of Bibliographic References.'
$ws.Range("C15").Value = '1) Scientific Methodology. 2) Monography Project. 3) Research Methods. 4) Guidelines of a research project. 5) Research in Bibliographic Databases. 6) Organization of Bibliographic References.'

# Row 16
$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = '1 – Metodologia Cientifica: Concepção e definição. 2 – Monografia Cientifica: O que é um projeto de pesquisa. As etapas de um projeto de pesquisa. 3 – Métodos de Pesquisa utilizados na Engenharia Química. 4- Pontos essenciais de um projeto de TCC. 5 – Pontos essenciais de uma monografia de TCC. 6 – Normas para elaboração de do texto e das Referências Bibliográficas. 7 – Mecanismos de busca em Bases de Dados tais como Web of Science, Scopus, Science Direct, etc. 8 – O uso do aplicativo Mendeley como forma de organizar e formatar as referências bibliográficas na monografia.'
$ws.Range("C16").Value = '1 – Metodologia Cientifica: Concepção e definição. 2 – Monografia Cientifica: O que é um projeto de pesquisa. As etapas de um projeto de pesquisa. 3 – Métodos de Pesquisa utilizados na Engenharia Química. 4- Pontos essenciais de um projeto de TCC. 5 – Pontos essenciais de uma monografia de TCC. 6 – Normas para elaboração de do texto e das Referências Bibliográficas. 7 – Mecanismos de busca em Bases de Dados tais como Web of Science, Scopus, Science Direct, etc. 8 – O uso do aplicativo Mendeley como forma de organizar e formatar as referências bibliográficas na monografia.'

# Row 17
$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("B17").Value = '1 - Scientific Methodology: design and definition. 2 - Scientific Monography: What is a research project. The steps of a research project. 3 - Research Methods used in Chemical Engineering. 4 Key points of a Course Conclusion Paper project. 5 - Key points of a Course Conclusion Paper monography. 6 - Standards for preparation of the text and the references. 7 - Search engines in Databases such as Web of Science, Scopus, Science Direct, etc. 8 - The use of the Mendeley application as a way of organizing and formatting the bibliographic references in the monograph.'
$ws.Range("C17").Value = '1 - Scientific Methodology: design and definition. 2 - Scientific Monography: What is a research project. The steps of a research project. 3 - Research Methods used in Chemical Engineering. 4 Key points of a Course Conclusion Paper project. 5 - Key points of a Course Conclusion Paper monography. 6 - Standards for preparation of the text and the references. 7 - Search engines in Databases such as Web of Science, Scopus, Science Direct, etc. 8 - The use of the Mendeley application as a way of organizing and formatting the bibliographic references in the monograph.'

# Row 18
$ws.Range("A18").Value = 'Avaliação:'

# Row 19
$ws.Range("A19").Value = 'Método:'
$ws.Range("B19").Value = 'Preparo e apresentação do Projeto de Trabalho de Conclusão de Curso (TCC 1) a ser desenvolvido na disciplina de Trabalho de Conclusão de Curso II, conforme norma do Departamento de Engenharia Química.'
$ws.Range("C19").Value = 'Preparo e apresentação do Projeto de Trabalho de Conclusão de Curso (TCC 1) a ser desenvolvido na disciplina de Trabalho de Conclusão de Curso II, conforme norma do Departamento de Engenharia Química.'

# Row 20
$ws.Range("A20").Value = 'Critério:'
$ws.Range("B20").Value = 'O aluno deve entregar, através do sistema online disponibilizado, um arquivo em formato pdf contendo o seu projeto de TCC impreterivelmente até a data estabelecida pelo professor na primeira semana de aula. O aluno que não cumprir este prazo fica reprovado na disciplina, por obter nota zero na primeira avaliação. O trabalho é submetido a dois avaliadores, a nota da primeira avaliação será a média das duas avaliações, sendo igual ou superior a 5, o aluno está aprovado, sendo inferior a 5 e igual ou superior a 3, o aluno está de recuperação.'
$ws.Range("C20").Value = 'O aluno deve entregar, através do sistema online disponibilizado, um arquivo em formato pdf contendo o seu projeto de TCC impreterivelmente até a data estabelecida pelo professor na primeira semana de aula. O aluno que não cumprir este prazo fica reprovado na disciplina, por obter nota zero na primeira avaliação. O trabalho é submetido a dois avaliadores, a nota da primeira avaliação será a média das duas avaliações, sendo igual ou superior a 5, o aluno está aprovado, sendo inferior a 5 e igual ou superior a 3, o aluno está de recuperação.'

# Row 21
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Range("B21").Value = 'O aluno deverá reapresentar o seu projeto até a data estabelecida pelo professor. O projeto será reavaliado e obtendo nota igual ou superior a 5, está aprovado.'
$ws.Range("C21").Value = 'O aluno deverá reapresentar o seu projeto até a data estabelecida pelo professor. O projeto será reavaliado e obtendo nota igual ou superior a 5, está aprovado.'

# Row 22
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = 'NASCIMENTO, L. P. Elaboração de Projetos de Pesquisa, Cengage Learning, 2012.SANTOS, C. R. Trabalho de Conclusão de Curso – Guia de elaboração passo a passo, Cengage Learning, 2010.BOOTH, W.; COLOMB, G.; WILLIAMS, J. A arte da Pesquisa. 3 ed. Martins Fontes. São Paulo. 2005. GIL, A.C. Como elaborar projetos de pesquisa. 5ed. Atlas, São Paulo, 2010.'
$ws.Range("C22").Value = 'NASCIMENTO, L. P. Elaboração de Projetos de Pesquisa, Cengage Learning, 2012.SANTOS, C. R. Trabalho de Conclusão de Curso – Guia de elaboração passo a passo, Cengage Learning, 2010.BOOTH, W.; COLOMB, G.; WILLIAMS, J. A arte da Pesquisa. 3 ed. Martins Fontes. São Paulo. 2005. GIL, A.C. Como elaborar projetos de pesquisa. 5ed. Atlas, São Paulo, 2010.'
